$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Plain value edits (fuzz-bug corrections)
$ws.Range("F13").Value = 325268233
$ws.Range("F14").Value = 644875817
$ws.Range("F15").Value = 6235000000
$ws.Range("F16").Value = 53616441

# Replace hard-coded totals with live formulas (matching the other columns)
$ws.Range("F18").Formula = "=SUM(F12:F17)"
$ws.Range("F21").Formula = "=SUM(F18:F20)"

$excel.CalculateFullRebuild()
$wb.Save()
